$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1612.9155
$ws.Range("I15").Value = 1612.9155
$ws.Range("K15").Value = 4838.7465
$ws.Range("M15").Value = -4669.7465
$ws.Range("H46").Value = 3323.75
$ws.Range("J46").Value = 2999.5
$ws.Range("L46").Value = 8998.5
$ws.Range("N46").Value = -9236.5
$ws.Range("H60").Value = 3323.75
$ws.Range("J60").Value = 2999.5
$ws.Range("L60").Value = 8998.5
$ws.Range("N60").Value = -9966.5
$ws.Range("H62").Value = 7112.161
$ws.Range("I62").Value = 6323.353
$ws.Range("K62").Value = 6323.353
$ws.Range("M62").Value = -5699.353
$ws.Range("H65").Value = 7112.161
$ws.Range("I65").Value = 6323.353
$ws.Range("K65").Value = 31616.765
$ws.Range("M65").Value = -28496.765
$ws.Range("H80").Value = 2181.205
$ws.Range("J80").Value = 3156.4736
$ws.Range("L80").Value = 9469.4208
$ws.Range("N80").Value = -11465.4208
$ws.Range("H83").Value = 2181.205
$ws.Range("J83").Value = 3156.4736
$ws.Range("L83").Value = 28408.2624
$ws.Range("N83").Value = -38392.2624
$ws.Range("H86").Value = 5850
$ws.Range("I86").Value = 4800.2
$ws.Range("J86").Value = 6433.222
$ws.Range("K86").Value = 4800.2
$ws.Range("L86").Value = 6433.222
$ws.Range("M86").Value = -3677.2
$ws.Range("N86").Value = -8679.222
$ws.Range("H88").Value = 191009.5
$ws.Range("I88").Value = 3000000
$ws.Range("J88").Value = 3743.4666
$ws.Range("K88").Value = 3000000
$ws.Range("L88").Value = 3743.4666
$ws.Range("M88").Value = -2999594
$ws.Range("N88").Value = -4555.4666
$ws.Range("H89").Value = 5850
$ws.Range("I89").Value = 4800.2
$ws.Range("J89").Value = 6433.222
$ws.Range("K89").Value = 24001
$ws.Range("L89").Value = 32166.11
$ws.Range("M89").Value = -18385
$ws.Range("N89").Value = -43398.11
$ws.Range("H91").Value = 191009.5
$ws.Range("I91").Value = 3000000
$ws.Range("J91").Value = 3743.4666
$ws.Range("K91").Value = 3000000
$ws.Range("L91").Value = 3743.4666
$ws.Range("M91").Value = -2998596
$ws.Range("N91").Value = -6551.4666
$ws.Range("H107").Value = 419.33334
$ws.Range("I107").Value = 410.57144
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 410.57144
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1509.42856
$ws.Range("N107").Value = -4290
$ws.Range("H112").Value = 2057.3809
$ws.Range("J112").Value = 2089
$ws.Range("L112").Value = 6267
$ws.Range("N112").Value = -8483
$ws.Range("H138").Value = 3378.7727
$ws.Range("J138").Value = 3693.1226
$ws.Range("L138").Value = 11079.3678
$ws.Range("N138").Value = -21359.3678

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 20836082
$ws.Range("I74").Value = 27779004
$ws.Range("J74").Value = 7316
$ws.Range("K74").Value = 27779004
$ws.Range("L74").Value = 7316
$ws.Range("M74").Value = -27778130
$ws.Range("N74").Value = -9064
$ws.Range("H77").Value = 20836082
$ws.Range("I77").Value = 27779004
$ws.Range("J77").Value = 7316
$ws.Range("K77").Value = 138895020
$ws.Range("L77").Value = 36580
$ws.Range("M77").Value = -138890652
$ws.Range("N77").Value = -45316
$ws.Range("H97").Value = 876.9
$ws.Range("I97").Value = 981.1875
$ws.Range("J97").Value = 459.75
$ws.Range("K97").Value = 981.1875
$ws.Range("L97").Value = 459.75
$ws.Range("M97").Value = -485.1875
$ws.Range("N97").Value = -1451.75
$ws.Range("H109").Value = 96250
$ws.Range("J109").Value = 96250
$ws.Range("L109").Value = 96250
$ws.Range("N109").Value = -99024
$ws.Range("H133").Value = 74874
$ws.Range("J133").Value = 74874
$ws.Range("L133").Value = 74874
$ws.Range("N133").Value = -79934
$ws.Range("H135").Value = 59332.332
$ws.Range("J135").Value = 59332.332
$ws.Range("L135").Value = 59332.332
$ws.Range("N135").Value = -69472.33199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5163.8076
$ws.Range("I20").Value = 4608.8423
$ws.Range("K20").Value = 4608.8423
$ws.Range("M20").Value = -4361.8423
$ws.Range("H41").Value = 200000
$ws.Range("J41").Value = 200000
$ws.Range("L41").Value = 200000
$ws.Range("N41").Value = -200776
$ws.Range("H94").Value = 1629.875
$ws.Range("I94").Value = 1405.2
$ws.Range("K94").Value = 1405.2
$ws.Range("M94").Value = -954.2
$ws.Range("H105").Value = 15116.24
$ws.Range("I105").Value = 13804.765
$ws.Range("K105").Value = 13804.765
$ws.Range("M105").Value = -12057.765

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9633.799999999999
$ws.Range("I62").Value = 3761.7144
$ws.Range("J62").Value = 23335.334
$ws.Range("K62").Value = 3761.7144
$ws.Range("L62").Value = 23335.334
$ws.Range("M62").Value = -3137.7144
$ws.Range("N62").Value = -24583.334
$ws.Range("H65").Value = 9633.799999999999
$ws.Range("I65").Value = 3761.7144
$ws.Range("J65").Value = 23335.334
$ws.Range("K65").Value = 18808.572
$ws.Range("L65").Value = 116676.67
$ws.Range("M65").Value = -15688.572
$ws.Range("N65").Value = -122916.67
$ws.Range("H99").Value = 2305.4443
$ws.Range("J99").Value = 2687.25
$ws.Range("L99").Value = 2687.25
$ws.Range("N99").Value = -5683.25
$ws.Range("H112").Value = 63157
$ws.Range("J112").Value = 63157
$ws.Range("L112").Value = 63157
$ws.Range("N112").Value = -66111
$ws.Range("H126").Value = 2305.4443
$ws.Range("J126").Value = 2687.25
$ws.Range("L126").Value = 8061.75
$ws.Range("N126").Value = -13001.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1315.3334
$ws.Range("J113").Value = 1470.375
$ws.Range("L113").Value = 4411.125
$ws.Range("N113").Value = -8751.125
$ws.Range("H136").Value = 3206
$ws.Range("I136").Value = 2757.5
$ws.Range("K136").Value = 8272.5
$ws.Range("M136").Value = -3172.5
$ws.Range("H137").Value = 2853
$ws.Range("I137").Value = 1566.2
$ws.Range("J137").Value = 3657.25
$ws.Range("K137").Value = 4698.6
$ws.Range("L137").Value = 10971.75
$ws.Range("M137").Value = 401.3999999999996
$ws.Range("N137").Value = -21171.75
$ws.Range("H139").Value = 6046.1113
$ws.Range("I139").Value = 5375
$ws.Range("K139").Value = 16125
$ws.Range("M139").Value = -10985

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 316969.38
$ws.Range("I80").Value = 502249.9
$ws.Range("K80").Value = 502249.9
$ws.Range("M80").Value = -501251.9
$ws.Range("H83").Value = 316969.38
$ws.Range("I83").Value = 502249.9
$ws.Range("K83").Value = 2511249.5
$ws.Range("M83").Value = -2506257.5
$ws.Range("H97").Value = 2878.5
$ws.Range("I97").Value = 2717.875
$ws.Range("K97").Value = 2717.875
$ws.Range("M97").Value = -2221.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10208.5
$ws.Range("J46").Value = 10714.571
$ws.Range("L46").Value = 10714.571
$ws.Range("N46").Value = -11090.571
$ws.Range("H132").Value = 4067
$ws.Range("I132").Value = 4028.4707
$ws.Range("K132").Value = 12085.4121
$ws.Range("M132").Value = -9555.4121

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 50000
$ws.Range("K87").Value = 50000
$ws.Range("M87").Value = -48752
$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 50000
$ws.Range("K90").Value = 150000
$ws.Range("M90").Value = -143760
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H112").Value = 53946.5
$ws.Range("J112").Value = 53946.5
$ws.Range("L112").Value = 53946.5
$ws.Range("N112").Value = -56900.5
$ws.Range("H132").Value = 2907.72
$ws.Range("I132").Value = 2273.7368
$ws.Range("K132").Value = 6821.2104
$ws.Range("M132").Value = -4291.2104
$ws.Range("H136").Value = 5110.533
$ws.Range("I136").Value = 977.2
$ws.Range("K136").Value = 2931.6
$ws.Range("M136").Value = -381.6000000000004
